# NA haver pull update
# Updates a set of data cells on the active sheet of the national accounts
# workbook with revised figures from the latest data pull.
#
# Every target cell already holds its number formatted as text (inline
# string) in the source workbook, so we force text ("@") number format
# before assigning the new value. This guarantees the literal string we
# assign (including trailing zeros like "71.0" or leading padding spaces
# like " 768678") round-trips exactly instead of being normalized as a
# numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "I2"  = "21157.1"
    "H4"  = "1600.1"
    "I4"  = "1687.2"
    "H6"  = "1369.2"
    "I6"  = "1429.6"
    "I7"  = "205.3"
    "I8"  = "499.4"
    "I9"  = "1317.3"
    "B12" = "63.5"
    "C12" = "68.5"
    "D12" = "68.7"
    "E12" = "69.8"
    "F12" = "71.0"
    "G12" = "62.7"
    "H12" = "54.2"
    "I12" = "79.1"
    "I16" = " 768678"
    "I18" = " 58400"
    "H20" = "609300"
    "I20" = "865600"
    "H21" = "63800"
    "I21" = "15000"
    "H22" = "73300"
    "I22" = "73300"
    "H23" = "140000"
    "I23" = "140000"
    "I25" = "1486802"
    "I26" = "1335322"
    "I28" = "835999"
    "I29" = "691267"
    "I30" = "1992677"
    "I31" = "2329677"
    "I34" = " 39200"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
